$wb = $excel.ActiveWorkbook

# === YDS sheet: append Week 16 game log numbers to the running lists ===
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "9 13 1 3 9 6 7 6 5 0 3 9 6 0 6 4 4 2 0 6 2 2 1 6 4 -1 12 7 11 5 1 6 -1 5 4 -1 2 -2 4 -2 4 5 4 18 -1 6 17 8 10 3 21 6 6 4 4 4 10 9 5 11 3 6 2 4 11 2 11 12 6 1 3 5 5 20 1 3 0 7 2 3 0 8 0 2 4 1 4 2 0 5 58 7 4 1 3 1 6 -1 9 0 5 10 12 2 9 2 10 8 9 8 5 2 1 4 3 4 8 6 7 0 5 3 1 2 4 24 1 4 9 3 12 1 2 2 5 -1 4 14 4 1 3 2 4 4 3 4 8 8 1 4 3 7 4 5 6 3 2 3 -4 -1 15 1 -2 3 0 4 11 7 -1 5 1 6 2 8 2 0 -2 3 5 5 1 2 0 66 3 5 12 4 5 -1 8 8 10 3 8 14 0 18 3 1 1 1 2 -1 1 7 10 -3 2 7 2 2 4 4 4 4 4 6 1 6 0 5 6 3 9 7 3 1 8 8 1 4 4 4 3 2 12 7 9 7 3 -1 10 2 10 4 4 2 1 3 1 -1 10 0 2 9 -1 7 1 1 -3 13 -1 1 2 -1 4 0 2 1 5 2 -4 3 5 2 2 -1 0 1 3 10 3 0 7 12 16 0 4 1 2 -1 8 8 5 1 1 9 3 -1 2 5 -1 10 3 0 9 3 1 1 6 6 6 0 8 0 2 4 3 2 4 7 1 6 3 9 7 2 3 6 1 0 26"
$wsYDS.Range("B3").Value = "13 33 15 8 3 -1 27 22 7 6 5 -3 7 15 41 7 9 6 9 5 18 0 13 9 22 18 4 4 24 7 13 25 -2 -8 6 5 19 0 6 6 13 6 0 8 14 11 20 7 9 12 13 5 12 11 6 12 2 12 9 13 6 9 8 17 7 12 10 7 7 52 7 4 27 5 20 -2 7 11 5 18 5 -2 25 7 22 13 7 13 7 18 58 9 4 8 4 5 12 19 7 8 24 1 19 4 6 8 9 11 20 10 28 11 29 11 11 21 7 14 10 7 10 20 12 6 1 7 5 17 9 10 13 7 7 10 3 9 7 0 7 6 3 6 11 9 8 11 7 4 5 9 10 9 5 8 11 4 6 12 -1 6 7 5 5 4 5 28 9 4 16 9 0 12 24 15 9 5 6 13 18 0 23 4 14 9 9 9 6 6 18 16 6 8 15 7 3 7 11 13 22 8 9 1 11 5 9 8 11 24 -1 21 -2 8 9 14 5 14 -5 21 7 11 26 20 3 8 3 28 -1 6 11 5 16 9 -2 12 9 7 17 10 9 8 14 3 7 5 7 3 12 3 31 7 14 6 5 5 9 9 31 41 6 6 3 2 30 11 12 21 9 8 3 6 5 10 11 2 11 7 5 5 7 5 5 11 1 4 12 9 9 8 17 10 11 1 7 5 10 25 19 11 13 22 18 17 8 7 6 6 4"
$wsYDS.Range("C2").Value = "6 6 5 7 -1 2 1 1 6 3 8 2 4 4 1 1 5 3 3 0 0 5 3 3 11 4 5 5 2 4 9 -4 29 4 5 0 3 3 0 3 0 5 9 3 2 7 0 1 -6 0 3 3 16 2 2 5 0 5 8 -1 3 14 2 0 3 6 5 2 5 1 1 2 0 4 7 6 6 1 9 4 2 4 4 4 0 8 11 4 0 0 -1 3 1 1 3 -3 2 2 2 7 1 0 4 3 9 4 11 5 7 0 1 4 0 1 8 1 2 4 2 4 4 4 6 -1 0 1 5 9 5 11 1 6 9 4 2 3 3 8 7 13 9 1 4 -2 0 4 15 0 9 5 4 11 10 -1 3 7 4 5 2 16 4 3 4 1 8 1 7 5 0 1 3 4 12 9 2 5 0 0 0 2 7 5 -5 4 3 -5 2 2 3 7 6 4 2 -2 -2 7 15 12 -3 1 22 4 -1 1 7 5 4 2 3 -1 34 0 8 3 7 20 5 8 4 4 9 3 1 3 6 4 0 6 0 -1 3 2 0 1 4 4 1 2 2 2 4 0 12 3 2 9 0 1 3 11 -1 25 3 5 -1 1 0 4 4 3 11 21 -2 12 0 5 2 4 4 -1 3 2 3 7 0 7 5 -2 4 19 15 1 11 7 1 27 4 12 4 2 1 2 7 0 6 1 2 6 2 5 2 7 6 0 5 9 2 2 9 4 9 -2 3 0 3 6 8 6 -1 6 7 5 0 25 2 3 2 2 6 3 17 0 5 4 4 6 8 2 5 3 4 3 2 -5 3 1 1 -1 2 4 5 6 1 1 12 0 2 -1 2 3 6 2 0 3 4 2 0 8 3 6 3 4 4 3 6 3 -1 0 2 5 4 2 1 0 8 3 3 52 1 6 3 14 3 3 1 4 5 0 14 9 4 12 -1 3 20 7 5 10 0 11 2 0 3 1 38 3 3 2"
$wsYDS.Range("C3").Value = "10 8 7 40 17 25 7 24 12 10 7 2 52 8 7 9 3 29 4 6 4 6 13 14 8 6 6 7 4 10 16 33 15 12 9 55 14 5 10 11 9 0 4 11 14 6 7 7 5 28 9 7 8 6 10 10 -2 12 5 6 34 8 7 9 18 33 13 1 36 3 7 19 7 33 -2 15 18 0 10 21 0 44 22 25 12 1 11 8 4 8 12 7 31 12 15 9 7 25 10 12 25 10 22 22 8 14 26 20 6 7 1 14 7 10 17 4 20 6 6 6 22 25 15 3 2 8 1 5 20 7 8 5 14 0 6 32 20 -1 12 18 2 12 9 8 5 9 27 1 5 8 5 9 5 8 16 7 28 5 4 8 5 28 5 9 3 7 1 16 0 4 9 9 14 1 7 2 9 11 13 7 9 5 22 8 12 15 3 5 -2 28 1 7 19 5 14 3 4 7 5 7 12 2 18 15 9 3 6 10 6 15 -5 1 -3 10 16 3 27 12 9 7 8 13 34 7 6 9 9 23 15 7 13 4 11 6 9 7 18 13 7 0 4 7 15 10 20 19 5 12 8 14 9 7 -1 27 4 35 5 9 8 7 7 7 4 19 8 7 43 29 8 3 19 22 2 5 5 3 8 1 9 6 11 15 1 16 9 12 5 4 8 11 4 10 24 9 5 13 7 12 10 5 6 3 7 22 7 7 1 7 6 12 11 33 2 6 18 3 43 4 19 2 4 4 6 4 22 2 7 12 5"

# === ST sheet: append Week 16 game log numbers to the running lists ===
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value = "65 66 70 58 64 62 63 64 62 63 66 60 64 61 61 57 63 62 63 64 64 64 63 58 55 63 67 60 68 60 52"
$wsST.Range("B5").Value = "24 23 19 24 18 20 21 45 33 28 23 27 19 26 21 24 22 26 19 65 16 25 98 0 23 20 2 34 31 29 19"
$wsST.Range("B6").Value = "18 2 18 27 24 16 24 21 23 27 5 20 15 13 25 22 21 28 19 20 25 23 28 29 26 19 24 0 16 20 16 23"
$wsST.Range("D3").Value = "55 52 42 49 41 58 44 63 56 64 55 44 47 56 36 43 43 57 40 42 38 45 27 37 43 49 51 54 41 53 48 33 46 42 48 46 42 41 54 59 51 52 51 39 44 42 44 68 53 31 55 40 49 32 47 43 54 54 40 54 40 49 59"
$wsST.Range("D4").Value = "0 8 0 0 0 9 0 4 0 25 27 0 0 0 0 12 0 -1 0 0 0 0 0 0 0 0 0 20 0 0 0 0 0 0 0 0 11 0 0 0 1 0 0 0 0 0 3 0 0 0 0 0 2 0 0 1 10 16 3 0 0 14 12"
$wsST.Range("D5").Value = "0 0 0 1 0 0 0 2 0 0 0 0 0 10 0 0 0 0 0 0 0 0 0 0 0 0 0 13 0 0 4 14 0 0 0 0 0 8 6 10 12 0 0 0 18 0 11 0 0 0 0 0 0 0 0 -2 48 3 0 4"

# === OFF sheet: update season totals (Home row 2, Road row 3) ===
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 148
$wsOFF.Range("D2").Value = 13
$wsOFF.Range("F2").Value = 53
$wsOFF.Range("G2").Value = 51
$wsOFF.Range("J2").Value = 25
$wsOFF.Range("N2").Value = 12
$wsOFF.Range("O2").Value = 30
$wsOFF.Range("P2").Value = 17
$wsOFF.Range("B3").Value = 25
$wsOFF.Range("C3").Value = 162
$wsOFF.Range("D3").Value = 2
$wsOFF.Range("E3").Value = 32
$wsOFF.Range("F3").Value = 84
$wsOFF.Range("G3").Value = 35
$wsOFF.Range("H3").Value = 27
$wsOFF.Range("I3").Value = 52
$wsOFF.Range("J3").Value = 48
$wsOFF.Range("L3").Value = 272
$wsOFF.Range("M3").Value = 160
$wsOFF.Range("Q3").Value = 448

# === DEF sheet: update season totals (Home row 2, Road row 3) ===
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 207
$wsDEF.Range("D2").Value = 12
$wsDEF.Range("F2").Value = 63
$wsDEF.Range("G2").Value = 59
$wsDEF.Range("J2").Value = 34
$wsDEF.Range("N2").Value = 12
$wsDEF.Range("O2").Value = 15
$wsDEF.Range("P2").Value = 12
$wsDEF.Range("C3").Value = 133
$wsDEF.Range("E3").Value = 30
$wsDEF.Range("F3").Value = 97
$wsDEF.Range("G3").Value = 37
$wsDEF.Range("H3").Value = 27
$wsDEF.Range("I3").Value = 56
$wsDEF.Range("J3").Value = 58
$wsDEF.Range("L3").Value = 215
$wsDEF.Range("M3").Value = 148
$wsDEF.Range("Q3").Value = 437

# === ST sheet: update season totals (row 2) ===
$wsST.Range("B2").Value = 50
$wsST.Range("D2").Value = 63
$wsST.Range("H2").Value = 6
$wsST.Range("J2").Value = 47
$wsST.Range("K2").Value = 41
$wsST.Range("L2").Value = 33
$wsST.Range("M2").Value = 24

# === TURNS sheet: update season totals (Road row 3) ===
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("D3").Value = 11

# === PEN sheet: update season totals ===
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 24
$wsPEN.Range("B3").Value = 17
